$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '56.828.66'
$ws.Range("E2").Value = '  +2.43%  '

# Row 3
$ws.Range("D3").Value = '3.003.74'
$ws.Range("E3").Value = '  +1.87%  '

# Row 4
$ws.Range("E4").Value = '  +0.13%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '512.97'
$ws.Range("E5").Value = '  +4.69%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.81'
$ws.Range("E6").Value = '  +4.81%  '

# Row 7
$ws.Range("E7").Value = '  -0.12%  '

# Row 8
$ws.Range("E8").Value = '  +3.28%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.54'
$ws.Range("E9").Value = '  +5.99%  '

# Row 10
$ws.Range("E10").Value = '  +7.17%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.357'
$ws.Range("E11").Value = '  +3.14%  '

# Row 12
$ws.Range("E12").Value = '  +2.61%  '

# Row 13
$ws.Range("D13").Value = '3.517.01'
$ws.Range("E13").Value = '  +2.04%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.66'
$ws.Range("E14").Value = '  +4.36%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000157'
$ws.Range("E15").Value = '  +11.81%  '

# Row 16
$ws.Range("D16").Value = '56.880.31'
$ws.Range("E16").Value = '  +2.58%  '

# Row 17
$ws.Range("D17").Value = '3.004.75'
$ws.Range("E17").Value = '  +1.87%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.94'
$ws.Range("E18").Value = '  +5.28%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.57'
$ws.Range("E19").Value = '  +3.75%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.86'
$ws.Range("E20").Value = '  +4.38%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.83'
$ws.Range("E21").Value = '  +3.21%  '

# Row 22
$ws.Range("E22").Value = '  +0.04%  '

# Row 23
$ws.Range("E23").Value = '  +5.57%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.38'
$ws.Range("E24").Value = '  +5.78%  '

# Row 25
$ws.Range("E25").Value = '  +7.15%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.93%  '

# Row 27
$ws.Range("D27").Value = '0.0₃0918'
$ws.Range("E27").Value = '  +9.19%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.63'
$ws.Range("E28").Value = '  +2.88%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.09'
$ws.Range("E29").Value = '  +8.40%  '

# Row 30
$ws.Range("E30").Value = '  +6.30%  '

# Row 31
$ws.Range("E31").Value = '  +6.73%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.61'
$ws.Range("E32").Value = '  +6.03%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '157.67'
$ws.Range("E33").Value = '  +5.98%  '

# Row 34
$ws.Range("E34").Value = '  +4.83%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.71'
$ws.Range("E35").Value = '  +0.93%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.28'
$ws.Range("E36").Value = '  -1.82%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0679'
$ws.Range("E37").Value = '  +3.85%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '24.09'
$ws.Range("E38").Value = '  +4.42%  '

# Row 39
$ws.Range("D39").Value = '3.035.69'
$ws.Range("E39").Value = '  +2.07%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.14'
$ws.Range("E40").Value = '  +2.64%  '

# Row 41
$ws.Range("E41").Value = '  +0.16%  '

# Row 42
$ws.Range("D42").Value = '2.286.23'
$ws.Range("E42").Value = '  +8.48%  '

# Row 43
$ws.Range("E43").Value = '  +3.64%  '

# Row 44
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.42'
$ws.Range("E44").Value = '  +3.13%  '

# Row 45
$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.68'
$ws.Range("E45").Value = '  +4.43%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +1.20%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.97'
$ws.Range("E47").Value = '  +9.36%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0240'
$ws.Range("E48").Value = '  +2.65%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.88'
$ws.Range("E49").Value = '  +6.36%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.20'
$ws.Range("E50").Value = '  -0.09%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0876'
$ws.Range("E51").Value = '  +4.50%  '
